$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.883.66'
$ws.Range("E2").Value = '  +1.92%  '
$ws.Range("D3").Value = '3.466.33'
$ws.Range("E3").Value = '  +1.80%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '161.49'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.78%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").Value = '3.468.25'
$ws.Range("E8").Value = '  +1.42%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.582'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +8.67%  '
$ws.Range("E10").Value = '  -2.53%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.125'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.47%  '
$ws.Range("E12").Value = '  +1.42%  '
$ws.Range("D13").Value = '4.071.84'
$ws.Range("E13").Value = '  +1.99%  '
$ws.Range("E14").Value = '  -2.78%  '
$ws.Range("E15").Value = '  +4.99%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.58'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.20%  '
$ws.Range("D17").Value = '64.911.07'
$ws.Range("E17").Value = '  +1.88%  '
$ws.Range("D18").Value = '3.478.59'
$ws.Range("E18").Value = '  +2.74%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.38'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.13%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.32'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '390.25'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.18'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.548'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.90'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.45%  '
$ws.Range("E25").Value = '  +0.68%  '
$ws.Range("E26").Value = '  +16.38%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.54'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("E28").Value = '  +0.16%  '
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.20'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +7.87%  '
$ws.Range("E31").Value = '  +7.83%  '
$ws.Range("E32").Value = '  -0.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.69'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.71%  '
$ws.Range("E34").Value = '  -1.15%  '
$ws.Range("E35").Value = '  +0.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.07'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.16%  '
$ws.Range("E37").Value = '  +0.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '161.89'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.32%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.90'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.79%  '
$ws.Range("D40").Value = '3.016.94'
$ws.Range("E40").Value = '  +3.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0768'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.98%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.39'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.56'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.36%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.81'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.10%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0317'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.18%  '
$ws.Range("E46").Value = '  +1.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.25'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +7.12%  '
$ws.Range("E48").Value = '  +1.94%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.880'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.32%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.18'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +10.29%  '
$ws.Range("E51").Value = '  +4.08%  '
